$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns to refreshed cryptos data.
# Price cells are forced as text (leading apostrophe) and the original
# cell style is restored afterwards so no stray style/number-format is introduced.

$style = $ws.Range("D2").Style
$ws.Range("D2").Value = "'27.625.81"
$ws.Range("D2").Style = $style
$ws.Range("E2").Value = '  -0.09%  '
$style = $ws.Range("D3").Style
$ws.Range("D3").Value = "'1.633.06"
$ws.Range("D3").Style = $style
$ws.Range("E3").Value = '  -0.29%  '
$ws.Range("E4").Value = '  +0.05%  '
$style = $ws.Range("D5").Style
$ws.Range("D5").Value = "'212.15"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = '  -0.20%  '
$style = $ws.Range("D6").Style
$ws.Range("D6").Value = "'0.523"
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = '  -0.10%  '
$ws.Range("E7").Value = '  +0.05%  '
$style = $ws.Range("D8").Style
$ws.Range("D8").Value = "'23.32"
$ws.Range("D8").Style = $style
$ws.Range("E8").Value = '  +1.40%  '
$ws.Range("E9").Value = '  +2.54%  '
$ws.Range("E10").Value = '  +0.23%  '
$style = $ws.Range("D11").Style
$ws.Range("D11").Value = "'0.0866"
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = '  -2.88%  '
$style = $ws.Range("D12").Style
$ws.Range("D12").Value = "'1.865.82"
$ws.Range("D12").Style = $style
$style = $ws.Range("D13").Style
$ws.Range("D13").Value = "'1.637.17"
$ws.Range("D13").Style = $style
$ws.Range("E13").Value = '  +0.04%  '
$ws.Range("E14").Value = '  +0.29%  '
$style = $ws.Range("D15").Style
$ws.Range("D15").Value = "'0.552"
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = '  -0.96%  '
$style = $ws.Range("D16").Style
$ws.Range("D16").Value = "'65.22"
$ws.Range("D16").Style = $style
$ws.Range("E16").Value = '  +1.01%  '
$style = $ws.Range("D17").Style
$ws.Range("D17").Value = "'27.618.19"
$ws.Range("D17").Style = $style
$ws.Range("E17").Value = '  -0.08%  '
$style = $ws.Range("D18").Style
$ws.Range("D18").Value = "'230.65"
$ws.Range("D18").Style = $style
$ws.Range("E18").Value = '  +0.31%  '
$ws.Range("E19").Value = '  -0.25%  '
$ws.Range("E20").Value = '  -2.21%  '
$ws.Range("E21").Value = '  -0.04%  '
$style = $ws.Range("D22").Style
$ws.Range("D22").Value = "'10.60"
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = '  +5.04%  '
$ws.Range("E23").Value = '  +1.38%  '
$style = $ws.Range("D24").Style
$ws.Range("D24").Value = "'2.11"
$ws.Range("D24").Style = $style
$ws.Range("E24").Value = '  +5.86%  '
$style = $ws.Range("D25").Style
$ws.Range("D25").Value = "'149.30"
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = '  -0.74%  '
$ws.Range("E26").Value = '  -1.02%  '
$ws.Range("E27").Value = '  -0.09%  '
$style = $ws.Range("D28").Style
$ws.Range("D28").Value = "'15.54"
$ws.Range("D28").Style = $style
$ws.Range("E28").Value = '  -0.41%  '
$ws.Range("E29").Value = '  +0.04%  '
$ws.Range("E30").Value = '  +0.10%  '
$ws.Range("E31").Value = '  -0.53%  '
$ws.Range("E32").Value = '  -0.53%  '
$style = $ws.Range("D33").Style
$ws.Range("D33").Value = "'1.482.23"
$ws.Range("D33").Style = $style
$ws.Range("E33").Value = '  +1.71%  '
$style = $ws.Range("D34").Style
$ws.Range("D34").Value = "'3.10"
$ws.Range("D34").Style = $style
$ws.Range("E34").Value = '  -0.31%  '
$ws.Range("E35").Value = '  -1.91%  '
$ws.Range("E36").Value = '  -1.18%  '
$style = $ws.Range("D37").Style
$ws.Range("D37").Value = "'0.937"
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = '  +3.83%  '
$ws.Range("E38").Value = '  -0.51%  '
$style = $ws.Range("D39").Style
$ws.Range("D39").Value = "'0.879"
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = '  +0.40%  '
$ws.Range("E40").Value = '  +0.21%  '
$ws.Range("E41").Value = '  +2.60%  '
$ws.Range("E42").Value = '  +0.01%  '
$style = $ws.Range("D43").Style
$ws.Range("D43").Value = "'67.94"
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = '  -3.08%  '
$ws.Range("E44").Value = '  +0.26%  '
$style = $ws.Range("D45").Style
$ws.Range("D45").Value = "'2.20"
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = '  -1.19%  '
$style = $ws.Range("D46").Style
$ws.Range("D46").Value = "'5.34"
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = '  -4.78%  '
$style = $ws.Range("D47").Style
$ws.Range("D47").Value = "'1.774.61"
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = '  -0.26%  '
$style = $ws.Range("D48").Style
$ws.Range("D48").Value = "'1.74"
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = '  +1.17%  '
$ws.Range("E49").Value = '  +1.45%  '
$ws.Range("E50").Value = '  -1.69%  '
$ws.Range("E51").Value = '  +0.70%  '
